$wb = $excel.ActiveWorkbook

# Sheet "Metadata": update URL, Version, Date, Publisher
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/hipaa-benefit-status"
$ws1.Range("B3").Value = "8.0.0"
$ws1.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$ws1.Range("B9").Value = "LinuxForHealth Team"

# Sheet "Include from HIPAA Benefit St": update System URI
$ws2 = $wb.Worksheets.Item("Include from HIPAA Benefit St")
$ws2.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/hipaa-benefit-status"
